$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 5 (ONT) before removing it - mirrors the user selecting the
# whole row in the UI prior to deleting it.
$ws.Rows("5:5").Select()

# Remove the ONT row entirely (row 5) - all following rows shift up by one.
$ws.Rows(5).Delete()

# Remove the MPL row entirely as well.
$mpl = $ws.Cells.Find("MPL")
if ($mpl -ne $null) {
    $mpl.EntireRow.Delete()
}

# ONDO's purchase numbers were corrected.
$ondo = $ws.Cells.Find("ONDO")
$ondoRow = $ondo.Row
$ws.Cells.Item($ondoRow, 2).Value = 3073
$ws.Cells.Item($ondoRow, 3).Value = 1.04

# Append the new AIOZ position as the last row of the table.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "AIOZ"
$ws.Cells.Item($newRow, 2).Value = 4621
$ws.Cells.Item($newRow, 3).Value = 0.5
$ws.Cells.Item($newRow, 4).Value = "KI"
